$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113, pushing the existing rows 113:175 down to 114:176.
$ws.Rows("113:113").Insert()

# Populate the newly inserted row 113 with the new price-report record.
$ws.Range("A113").Value = 8
$ws.Range("B113").Value = "Terminal La Palmera de La Serena"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44452
$ws.Range("E113").Value = 4
$ws.Range("F113").Value = 100112032
$ws.Range("G113").Value = "Zapallo italiano"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 13500
$ws.Range("L113").Value = 14000
$ws.Range("M113").Value = 13750
$ws.Range("N113").Value = "$/caja 50 unidades"
$ws.Range("O113").Value = "Región de Arica y Parinacota"
$ws.Range("P113").Value = 275
$ws.Range("Q113").Value = 50
$ws.Range("R113").Value = "Hortaliza"

# Match the date cell format used by the rest of column D (style index 2).
$ws.Range("D113").NumberFormat = $ws.Range("D114").NumberFormat
